$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Rename the affected labels (this causes the shared-string table to be
# rewritten: old entries "num", "Teil:Teil", "s.qlang" are replaced by
# new entries "Teil:Artikel", "num2", "S.qlang")
$ws.Range("A1").Value = "Teil:Artikel"
$ws.Range("A2").Value = "num2"
$ws.Range("D2").Value = "S.qlang"

# Move the active selection from D12 to D2
$ws.Range("D2").Select()

$wb.Save()
